$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $range = $ws.Range($cell)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue "D2" "27.658.00"
Set-TextValue "E2" "  +1.20%  "
Set-TextValue "D3" "1.643.20"
Set-TextValue "E3" "  -0.58%  "
Set-TextValue "E4" "  +0.04%  "
Set-TextValue "D5" "213.02"
Set-TextValue "E5" "  -0.11%  "
Set-TextValue "D6" "0.530"
Set-TextValue "E6" "  +3.19%  "
Set-TextValue "E7" "  +0.05%  "
Set-TextValue "D8" "23.00"
Set-TextValue "E8" "  -2.23%  "
Set-TextValue "E9" "  -1.20%  "
Set-TextValue "E10" "  -0.39%  "
Set-TextValue "D11" "0.0890"
Set-TextValue "D12" "1.879.64"
Set-TextValue "E12" "  -0.45%  "
Set-TextValue "D13" "1.643.04"
Set-TextValue "E13" "  -0.58%  "
Set-TextValue "E14" "  -0.68%  "
Set-TextValue "E15" "  -1.20%  "
Set-TextValue "D16" "64.13"
Set-TextValue "E16" "  -2.22%  "
Set-TextValue "D17" "27.585.37"
Set-TextValue "E17" "  +0.94%  "
Set-TextValue "D18" "229.16"
Set-TextValue "E18" "  -0.71%  "
Set-TextValue "E19" "  -0.51%  "
Set-TextValue "D20" "7.63"
Set-TextValue "E20" "  +2.28%  "
Set-TextValue "E21" "  +0.10%  "
Set-TextValue "E22" "  -1.31%  "
Set-TextValue "D23" "10.05"
Set-TextValue "E23" "  +7.54%  "
Set-TextValue "E24" "  -3.05%  "
Set-TextValue "D25" "149.18"
Set-TextValue "E25" "  +1.55%  "
Set-TextValue "D26" "6.95"
Set-TextValue "E26" "  -2.85%  "
Set-TextValue "D27" "0.113"
Set-TextValue "E27" "  +1.00%  "
Set-TextValue "E28" "  -0.03%  "
Set-TextValue "E29" "  -1.50%  "
Set-TextValue "E30" "  -0.35%  "
Set-TextValue "D31" "0.0484"
Set-TextValue "E31" "  -2.58%  "
Set-TextValue "D32" "3.29"
Set-TextValue "E32" "  -0.12%  "
Set-TextValue "D33" "3.17"
Set-TextValue "E33" "  +2.11%  "
Set-TextValue "D34" "1.438.02"
Set-TextValue "E34" "  -1.05%  "
Set-TextValue "D35" "1.58"
Set-TextValue "E35" "  +2.58%  "
Set-TextValue "D36" "2.35"
Set-TextValue "E36" "  -1.30%  "
Set-TextValue "E37" "  +0.19%  "
Set-TextValue "E38" "  -2.86%  "
Set-TextValue "E39" "  -1.13%  "
Set-TextValue "D40" "0.898"
Set-TextValue "E40" "  +14.52%  "
Set-TextValue "E41" "  -1.60%  "
Set-TextValue "D42" "1.00"
Set-TextValue "E42" "  +0.14%  "
Set-TextValue "D43" "5.70"
Set-TextValue "E43" "  +4.28%  "
Set-TextValue "B44" "MXToken"
Set-TextValue "C44" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D44" "2.26"
Set-TextValue "E44" "  +1.83%  "
Set-TextValue "B45" "mCoin"
Set-TextValue "C45" "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
Set-TextValue "D45" "2.47"
Set-TextValue "E45" "  -0.67%  "
Set-TextValue "D46" "65.30"
Set-TextValue "E46" "  +0.63%  "
Set-TextValue "D47" "1.787.08"
Set-TextValue "E47" "  -0.50%  "
Set-TextValue "E48" "  -1.36%  "
Set-TextValue "D49" "86.42"
Set-TextValue "E49" "  -1.87%  "
Set-TextValue "D50" "0.0₆0102"
Set-TextValue "E50" "  -3.97%  "
Set-TextValue "D51" "0.0986"
Set-TextValue "E51" "  -2.57%  "
